# Adds the "Interim Degrees" program row (row 52) to the programs sheet,
# based on Obergruber & Zierow (2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row data -----------------------------------------------------
# Column order of entry matters because it controls the order in which
# new strings are appended to xl/sharedStrings.xml (A, B, then J, I, F).

$ws.Range("A52").Value = "interimDegrees"
$ws.Range("B52").Value = "Interim Degrees"
$ws.Range("C52").Value = 1978
$ws.Range("D52").Value = "Education"
$ws.Range("E52").Value = 16

$url = "https://ideas.repec.org/a/eee/ecoedu/v75y2020ics0272775718307179.html"
$ws.Range("J52").Value = $url
$ws.Hyperlinks.Add($ws.Range("J52"), $url) | Out-Null
$ws.Range("J52").Style = $ws.Range("J51").Style

$ws.Range("I52").Value = "Obergruber & Zierow (2020)"
$ws.Range("I52").WrapText = $true

$description = "All German federal states introduced interim degrees, which award the degree of the middle track ""Realschule"" to everyone who completes the the 10th grade of the highest track of the German school system ""Gymnasium"". The idea behind this reform was provide a fallback option for students not being able to complete the high track. Obergruber & Zierow (2020) find that this reform incentivized students to stay on the highest track."
$ws.Range("F52").Value = $description
$ws.Range("F52").WrapText = $true

# Row height matches the other multi-line rows (135pt, as in the diff).
$ws.Rows.Item(52).RowHeight = 135

# --- Update selection / active cell to the new row --------------------
$ws.Activate()
$excel.Goto($ws.Range("A52"), $true) | Out-Null
